# ----------------------------------------------------------------------------
# Reproduces the commit "Add .py file + example to keep the formatting":
#   - Renames Sheet1 -> "United States Of America"
#   - Gives the sheet tab a color
#   - Zooms the sheet view to 160%
#   - Sets explicit (best-fit-like) column widths for columns A:O
#   - Styles the header row (A1:O1): bold white font, blue fill, thin border,
#     centered / top-aligned text
#   - Highlights the "Discounts" column (H2:H141) with a light grey fill
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet -----------------------------------------------
$ws.Name = "United States Of America"

# --- Tab color (theme accent5, ~ #5B9BD5) --------------------------------
$ws.Tab.Color = 13998939

# --- Zoom the sheet view to 160% -----------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 160

# --- Column widths (A:O) --------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.833333333333332
$ws.Columns.Item(2).ColumnWidth = 22.0
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 8.833333333333332
$ws.Columns.Item(6).ColumnWidth = 17.666666666666668
$ws.Columns.Item(7).ColumnWidth = 8.5
$ws.Columns.Item(8).ColumnWidth = 9.833333333333332
$ws.Columns.Item(9).ColumnWidth = 8.5
$ws.Columns.Item(10).ColumnWidth = 9.166666666666666
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 9.166666666666666
$ws.Columns.Item(13).ColumnWidth = 13.666666666666666
$ws.Columns.Item(14).ColumnWidth = 11.5
$ws.Columns.Item(15).ColumnWidth = 4.166666666666667

# --- Header row style (A1:O1) ---------------------------------------------
# Build the full look on a scratch cell first (so every property lands in a
# single paste operation and only one new font/fill/border/style gets added),
# then copy the formatting onto the header row and clean up the scratch cell.
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Font.Color = 16777215
$scratch.Interior.Color = 9524736
$scratch.Borders.LineStyle = 1
$scratch.Borders.Weight = 2
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Copy()
$header = $ws.Range("A1:O1")
$header.PasteSpecial(-4122)
$scratch.Clear()

# --- Highlight the Discounts column (H2:H141) ------------------------------
$scratch2 = $ws.Range("Z1")
$scratch2.Interior.Color = 14277081
$scratch2.Copy()
$hcol = $ws.Range("H2:H141")
$hcol.PasteSpecial(-4122)
$scratch2.Clear()
